# Updated cryptos list on Thu Nov 28 23:15:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Coin = $null; Link = $null; Price = "95.649.66";  Volume = "-0.35%" },
    @{ Row = 3;  Coin = $null; Link = $null; Price = "3.583.37";   Volume = "-2.15%" },
    @{ Row = 4;  Coin = $null; Link = $null; Price = $null;        Volume = "-0.13%" },
    @{ Row = 5;  Coin = $null; Link = $null; Price = "237.20";     Volume = "-1.50%" },
    @{ Row = 6;  Coin = $null; Link = $null; Price = "655.54";     Volume = "+2.09%" },
    @{ Row = 7;  Coin = $null; Link = $null; Price = "1.52";       Volume = "+2.98%" },
    @{ Row = 8;  Coin = $null; Link = $null; Price = "0.402";      Volume = "+0.38%" },
    @{ Row = 9;  Coin = $null; Link = $null; Price = $null;        Volume = "+0.04%" },
    @{ Row = 10; Coin = $null; Link = $null; Price = "1.03";       Volume = "+2.98%" },
    @{ Row = 11; Coin = $null; Link = $null; Price = "3.580.53";   Volume = "-2.10%" },
    @{ Row = 12; Coin = $null; Link = $null; Price = "43.02";      Volume = "-1.72%" },
    @{ Row = 13; Coin = $null; Link = $null; Price = $null;        Volume = "+0.98%" },
    @{ Row = 14; Coin = $null; Link = $null; Price = "6.48";       Volume = "+1.69%" },
    @{ Row = 15; Coin = $null; Link = $null; Price = "4.250.12";   Volume = "-2.58%" },
    @{ Row = 16; Coin = $null; Link = $null; Price = "95.519.05";  Volume = "-0.26%" },
    @{ Row = 17; Coin = $null; Link = $null; Price = $null;        Volume = "-0.33%" },
    @{ Row = 18; Coin = $null; Link = $null; Price = "3.584.47";   Volume = "-2.52%" },
    @{ Row = 19; Coin = $null; Link = $null; Price = "12.75";      Volume = "-5.65%" },
    @{ Row = 20; Coin = $null; Link = $null; Price = "7.76";       Volume = "-3.39%" },
    @{ Row = 21; Coin = $null; Link = $null; Price = "17.99";      Volume = "-4.11%" },
    @{ Row = 22; Coin = "Stellar"; Link = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; Price = "0.497"; Volume = "+2.88%" },
    @{ Row = 23; Coin = "SuiNetwork"; Link = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; Price = "3.46"; Volume = "+1.01%" },
    @{ Row = 24; Coin = $null; Link = $null; Price = "510.88";     Volume = "-1.62%" },
    @{ Row = 25; Coin = $null; Link = $null; Price = "7.08";       Volume = "+4.41%" },
    @{ Row = 26; Coin = $null; Link = $null; Price = $null;        Volume = "+0.02%" },
    @{ Row = 27; Coin = $null; Link = $null; Price = "95.91";      Volume = "-1.76%" },
    @{ Row = 28; Coin = $null; Link = $null; Price = "12.83";      Volume = "+2.01%" },
    @{ Row = 29; Coin = $null; Link = $null; Price = "3.776.88";   Volume = "-2.29%" },
    @{ Row = 30; Coin = $null; Link = $null; Price = "3.05";       Volume = "-4.52%" },
    @{ Row = 31; Coin = $null; Link = $null; Price = "0.149";      Volume = "+5.22%" },
    @{ Row = 32; Coin = $null; Link = $null; Price = "11.61";      Volume = "-0.34%" },
    @{ Row = 33; Coin = $null; Link = $null; Price = $null;        Volume = "+0.18%" },
    @{ Row = 34; Coin = $null; Link = $null; Price = "0.997";      Volume = "+0.30%" },
    @{ Row = 35; Coin = $null; Link = $null; Price = $null;        Volume = "-1.21%" },
    @{ Row = 36; Coin = $null; Link = $null; Price = "31.95";      Volume = "-2.29%" },
    @{ Row = 37; Coin = $null; Link = $null; Price = "1.70";       Volume = "+13.57%" },
    @{ Row = 38; Coin = "PolygonEcosystemToken"; Link = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; Price = "0.565"; Volume = "-2.14%" },
    @{ Row = 39; Coin = "RenderToken"; Link = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; Price = "8.63"; Volume = "+10.20%" },
    @{ Row = 40; Coin = $null; Link = $null; Price = "597.37";     Volume = "+5.78%" },
    @{ Row = 41; Coin = $null; Link = $null; Price = $null;        Volume = "-0.34%" },
    @{ Row = 42; Coin = $null; Link = $null; Price = $null;        Volume = "+0.07%" },
    @{ Row = 43; Coin = $null; Link = $null; Price = $null;        Volume = "+8.26%" },
    @{ Row = 44; Coin = $null; Link = $null; Price = "0.914";      Volume = "-4.39%" },
    @{ Row = 45; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "5.77"; Volume = "-0.36%" },
    @{ Row = 46; Coin = "Stacks"; Link = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; Price = "2.30"; Volume = "+3.68%" },
    @{ Row = 47; Coin = "EnergySwap"; Link = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Price = "34.35"; Volume = "+4.49%" },
    @{ Row = 48; Coin = "WhiteBITCoin"; Link = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; Price = "23.44"; Volume = "-1.30%" },
    @{ Row = 49; Coin = "VeChain"; Link = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; Price = "0.0418"; Volume = "-2.78%" },
    @{ Row = 50; Coin = $null; Link = $null; Price = "3.49";       Volume = "-0.01%" },
    @{ Row = 51; Coin = $null; Link = $null; Price = $null;        Volume = "-0.66%" }
)

foreach ($r in $rows) {
    if ($r.Coin -ne $null) {
        $ws.Cells.Item($r.Row, 2).Value = $r.Coin
    }
    if ($r.Link -ne $null) {
        $ws.Cells.Item($r.Row, 3).Value = $r.Link
    }
    if ($r.Price -ne $null) {
        $ws.Cells.Item($r.Row, 4).Value = "'" + $r.Price
    }
    if ($r.Volume -ne $null) {
        $ws.Cells.Item($r.Row, 5).Value = "  " + $r.Volume + "  "
    }
}
